$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.555.20'
$ws.Range('E2').Value = '  +0.90%  '

$ws.Range('D3').Value = '3.032.51'
$ws.Range('E3').Value = '  +2.43%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = "'384.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.13%  '

$ws.Range('D6').Value = "'102.76"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.44%  '

$ws.Range('E7').Value = '  -0.22%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').Value = "'0.587"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.71%  '

$ws.Range('D10').Value = "'36.81"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.48%  '

$ws.Range('E11').Value = '  +0.01%  '

$ws.Range('E12').Value = '  +0.62%  '

$ws.Range('D13').Value = '3.510.87'
$ws.Range('E13').Value = '  +2.58%  '

$ws.Range('D14').Value = "'18.68"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.98%  '

$ws.Range('D15').Value = "'7.76"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.19%  '

$ws.Range('D16').Value = '3.034.13'
$ws.Range('E16').Value = '  +2.25%  '

$ws.Range('D17').Value = "'0.972"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.79%  '

$ws.Range('D18').Value = "'10.59"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -11.54%  '

$ws.Range('D19').Value = '51.583.97'
$ws.Range('E19').Value = '  +0.84%  '

$ws.Range('E20').Value = '  -0.52%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0961'
$ws.Range('E21').Value = '  -0.09%  '

$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').Value = "'12.32"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.86%  '

$ws.Range('D23').Value = "'69.87"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.26%  '

$ws.Range('D24').Value = "'266.65"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.44%  '

$ws.Range('D25').Value = "'3.16"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.97%  '

$ws.Range('E26').Value = '  +5.06%  '

$ws.Range('D27').Value = "'7.44"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.52%  '

$ws.Range('D28').Value = "'0.173"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.30%  '

$ws.Range('D29').Value = "'26.30"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.55%  '

$ws.Range('E30').Value = '  +0.00%  '

$ws.Range('E31').Value = '  -2.40%  '

$ws.Range('D32').Value = "'10.27"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.59%  '

$ws.Range('E33').Value = '  -0.50%  '

$ws.Range('D34').Value = "'34.07"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.98%  '

$ws.Range('D35').Value = "'50.50"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.06%  '

$ws.Range('D36').Value = "'0.0448"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.40%  '

$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('D38').Value = "'3.37"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.44%  '

$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = "'16.97"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.19%  '

$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = "'0.283"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.37%  '

$ws.Range('E41').Value = '  +1.11%  '

$ws.Range('E42').Value = '  -0.34%  '

$ws.Range('E43').Value = '  +2.07%  '

$ws.Range('E44').Value = '  +0.67%  '

$ws.Range('E45').Value = '  +4.11%  '

$ws.Range('D46').Value = "'21.58"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.21%  '

$ws.Range('D47').Value = "'2.47"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.96%  '

$ws.Range('E48').Value = '  +3.82%  '

$ws.Range('D49').Value = '2.036.24'
$ws.Range('E49').Value = '  -0.54%  '

$ws.Range('D50').Value = '3.336.55'
$ws.Range('E50').Value = '  +2.68%  '

$ws.Range('D51').Value = "'0.207"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.47%  '
